# Add a new "metadata" worksheet after the existing "data" sheet and
# populate it with the panel query metadata, then refresh the F-column
# (time_taken) timestamps on the "data" sheet.

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item(1)

# --- Create the metadata sheet, positioned right after "data" ---
$ws = $wb.Worksheets.Add($null, $data)
$ws.Name = "metadata"

# Match header/index styling used on the "data" sheet (bold, centered,
# bordered header row; bordered index column).
$data.Range("B1:F1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)
$data.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

# --- Header row ---
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

# --- Data row ---
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Congenital muscular dystrophy"
$ws.Range("C2").Value = 207
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "2.16"
$ws.Range("E2").Value = "2021-09-02T17:03:33.301659Z"
$ws.Range("F2").Value = "2021-10-05 14:19:44.734652"
$ws.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/207/?format=json"

# Keep "data" as the active sheet (unchanged from before the edit).
$data.Activate()

# --- Refresh query timestamps on the "data" sheet (col F, rows 2-59) ---
$data.Range("F2").Value = "2021-10-05 14:19:44.738186"
$data.Range("F3").Value = "2021-10-05 14:19:44.738194"
$data.Range("F4").Value = "2021-10-05 14:19:44.738197"
$data.Range("F5").Value = "2021-10-05 14:19:44.738200"
$data.Range("F6").Value = "2021-10-05 14:19:44.738203"
$data.Range("F7").Value = "2021-10-05 14:19:44.738205"
$data.Range("F8").Value = "2021-10-05 14:19:44.738208"
$data.Range("F9").Value = "2021-10-05 14:19:44.738210"
$data.Range("F10").Value = "2021-10-05 14:19:44.738213"
$data.Range("F11").Value = "2021-10-05 14:19:44.738216"
$data.Range("F12").Value = "2021-10-05 14:19:44.738218"
$data.Range("F13").Value = "2021-10-05 14:19:44.738220"
$data.Range("F14").Value = "2021-10-05 14:19:44.738223"
$data.Range("F15").Value = "2021-10-05 14:19:44.738225"
$data.Range("F16").Value = "2021-10-05 14:19:44.738228"
$data.Range("F17").Value = "2021-10-05 14:19:44.738230"
$data.Range("F18").Value = "2021-10-05 14:19:44.738233"
$data.Range("F19").Value = "2021-10-05 14:19:44.738235"
$data.Range("F20").Value = "2021-10-05 14:19:44.738238"
$data.Range("F21").Value = "2021-10-05 14:19:44.738240"
$data.Range("F22").Value = "2021-10-05 14:19:44.738243"
$data.Range("F23").Value = "2021-10-05 14:19:44.738246"
$data.Range("F24").Value = "2021-10-05 14:19:44.738248"
$data.Range("F25").Value = "2021-10-05 14:19:44.738250"
$data.Range("F26").Value = "2021-10-05 14:19:44.738253"
$data.Range("F27").Value = "2021-10-05 14:19:44.738256"
$data.Range("F28").Value = "2021-10-05 14:19:44.738258"
$data.Range("F29").Value = "2021-10-05 14:19:44.738261"
$data.Range("F30").Value = "2021-10-05 14:19:44.738263"
$data.Range("F31").Value = "2021-10-05 14:19:44.738266"
$data.Range("F32").Value = "2021-10-05 14:19:44.738268"
$data.Range("F33").Value = "2021-10-05 14:19:44.738271"
$data.Range("F34").Value = "2021-10-05 14:19:44.738273"
$data.Range("F35").Value = "2021-10-05 14:19:44.738276"
$data.Range("F36").Value = "2021-10-05 14:19:44.738278"
$data.Range("F37").Value = "2021-10-05 14:19:44.738281"
$data.Range("F38").Value = "2021-10-05 14:19:44.738283"
$data.Range("F39").Value = "2021-10-05 14:19:44.738285"
$data.Range("F40").Value = "2021-10-05 14:19:44.738288"
$data.Range("F41").Value = "2021-10-05 14:19:44.738290"
$data.Range("F42").Value = "2021-10-05 14:19:44.738293"
$data.Range("F43").Value = "2021-10-05 14:19:44.738296"
$data.Range("F44").Value = "2021-10-05 14:19:44.738298"
$data.Range("F45").Value = "2021-10-05 14:19:44.738301"
$data.Range("F46").Value = "2021-10-05 14:19:44.738303"
$data.Range("F47").Value = "2021-10-05 14:19:44.738305"
$data.Range("F48").Value = "2021-10-05 14:19:44.738308"
$data.Range("F49").Value = "2021-10-05 14:19:44.738310"
$data.Range("F50").Value = "2021-10-05 14:19:44.738313"
$data.Range("F51").Value = "2021-10-05 14:19:44.738315"
$data.Range("F52").Value = "2021-10-05 14:19:44.738318"
$data.Range("F53").Value = "2021-10-05 14:19:44.738320"
$data.Range("F54").Value = "2021-10-05 14:19:44.738323"
$data.Range("F55").Value = "2021-10-05 14:19:44.738325"
$data.Range("F56").Value = "2021-10-05 14:19:44.738328"
$data.Range("F57").Value = "2021-10-05 14:19:44.738330"
$data.Range("F58").Value = "2021-10-05 14:19:44.738333"
$data.Range("F59").Value = "2021-10-05 14:19:44.738335"
